$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "Animation Designer (Januar 2021" + [char]0x2013 + "Heute)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Spark Animation: Animation Designer (Jan 2021 - Present)", 2)

$d.Content.Find.Execute(
    "Animation Designer (Juni 2018" + [char]0x2013 + "Dezember 2020)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Pixel Studio: Animations-Designer (Jun 2018 - Dez 2020)", 2)

$d.Content.Find.Execute(
    "Junior Animation Designer (September 2016" + [char]0x2013 + "Mai 2018)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Flash Animation: Junior Animation Designer (Sep 2016 - Mai 2018)", 2)

$d.Content.Find.Execute(
    "Master of Arts in Animation, voraussichtlicher Abschluss:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Master of Arts in Animation, erwartet Abschluss: Dez 2025", 2)
